$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.15163
$ws.Range("H2").Value = 9.454890000000001
$ws.Range("I2").Value = 0.0006291248881010851
$ws.Range("J2").Value = 0.0006291248881010851
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.307106666666667
$ws.Range("N2").Value = 3.92132
$ws.Range("O2").Value = 0.01256263154946851
$ws.Range("P2").Value = 0.01256263154946851
$ws.Range("Q2").Value = 4.119516583866667
$ws.Range("R2").Value = 37.07564925480001
$ws.Range("S2").Value = 0.000007903464167814538
$ws.Range("T2").Value = 0.000007903464167814538
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.15163
$ws.Range("H3").Value = 9.454890000000001
$ws.Range("I3").Value = 0.0006291248881010851
$ws.Range("J3").Value = 0.0006291248881010851
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("N3").Value = 240.678711
$ws.Range("O3").Value = 0.77105616682495
$ws.Range("P3").Value = 0.77105616682495
$ws.Range("Q3").Value = 252.84341531631
$ws.Range("R3").Value = 2275.59073784679
$ws.Range("S3").Value = 0.0004850906246733983
$ws.Range("T3").Value = 0.0004850906246733983
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.15163
$ws.Range("H4").Value = 9.454890000000001
$ws.Range("I4").Value = 0.0006291248881010851
$ws.Range("J4").Value = 0.0006291248881010851
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 22.51385866666667
$ws.Range("N4").Value = 67.54157600000001
$ws.Range("O4").Value = 0.2163812016255815
$ws.Range("P4").Value = 0.2163812016255815
$ws.Range("Q4").Value = 70.95535238962668
$ws.Range("R4").Value = 638.5981715066401
$ws.Range("S4").Value = 0.0001361307992598723
$ws.Range("T4").Value = 0.0001361307992598723
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4971.754394666666
$ws.Range("H5").Value = 14915.263184
$ws.Range("I5").Value = 0.9924561027819714
$ws.Range("J5").Value = 0.9924561027819713
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.307106666666667
$ws.Range("N5").Value = 3.92132
$ws.Range("O5").Value = 0.01256263154946851
$ws.Range("P5").Value = 0.01256263154946851
$ws.Range("Q5").Value = 6498.613314298097
$ws.Range("R5").Value = 58487.51982868288
$ws.Range("S5").Value = 0.01246786034827136
$ws.Range("T5").Value = 0.01246786034827136
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4971.754394666666
$ws.Range("H6").Value = 14915.263184
$ws.Range("I6").Value = 0.9924561027819714
$ws.Range("J6").Value = 0.9924561027819713
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("N6").Value = 240.678711
$ws.Range("O6").Value = 0.77105616682495
$ws.Range("P6").Value = 0.77105616682495
$ws.Range("Q6").Value = 398865.1463723194
$ws.Range("R6").Value = 3589786.317350875
$ws.Range("S6").Value = 0.7652393983530955
$ws.Range("T6").Value = 0.7652393983530954
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4971.754394666666
$ws.Range("H7").Value = 14915.263184
$ws.Range("I7").Value = 0.9924561027819714
$ws.Range("J7").Value = 0.9924561027819713
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 22.51385866666667
$ws.Range("N7").Value = 67.54157600000001
$ws.Range("O7").Value = 0.2163812016255815
$ws.Range("P7").Value = 0.2163812016255815
$ws.Range("Q7").Value = 111933.3757669042
$ws.Range("R7").Value = 1007400.381902138
$ws.Range("S7").Value = 0.2147488440806045
$ws.Range("T7").Value = 0.2147488440806045
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 34.63986933333333
$ws.Range("H8").Value = 103.919608
$ws.Range("I8").Value = 0.006914772329927541
$ws.Range("J8").Value = 0.006914772329927542
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.307106666666667
$ws.Range("N8").Value = 3.92132
$ws.Range("O8").Value = 0.01256263154946851
$ws.Range("P8").Value = 0.01256263154946851
$ws.Range("Q8").Value = 45.27800413806222
$ws.Range("R8").Value = 407.50203724256
$ws.Range("S8").Value = 0.00008686773702933961
$ws.Range("T8").Value = 0.00008686773702933963
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 34.63986933333333
$ws.Range("H9").Value = 103.919608
$ws.Range("I9").Value = 0.006914772329927541
$ws.Range("J9").Value = 0.006914772329927542
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("N9").Value = 240.678711
$ws.Range("O9").Value = 0.77105616682495
$ws.Range("P9").Value = 0.77105616682495
$ws.Range("Q9").Value = 2779.026366785031
$ws.Range("R9").Value = 25011.23730106528
$ws.Range("S9").Value = 0.005331677847181159
$ws.Range("T9").Value = 0.005331677847181159
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 34.63986933333333
$ws.Range("H10").Value = 103.919608
$ws.Range("I10").Value = 0.006914772329927541
$ws.Range("J10").Value = 0.006914772329927542
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 22.51385866666667
$ws.Range("N10").Value = 67.54157600000001
$ws.Range("O10").Value = 0.2163812016255815
$ws.Range("P10").Value = 0.2163812016255815
$ws.Range("Q10").Value = 779.8771224024675
$ws.Range("R10").Value = 7018.894101622209
$ws.Range("S10").Value = 0.001496226745717043
$ws.Range("T10").Value = 0.001496226745717043
